$wb = $excel.ActiveWorkbook

# --- "Comparison with GANs" sheet: drop the blank leading rows (1-8), which
# shifts the existing table up from rows 9:19 to rows 1:11, then refresh the
# handful of cells whose figures were updated. ---
$ws = $wb.Worksheets.Item("Comparison with GANs")
$ws.Range("A1:A8").EntireRow.Delete()

# Row 5 (was row 13)
$ws.Range("D5").Value = 2.92
$ws.Range("E5").Value = 2.79
$ws.Range("F5").Value = 1.13

# Row 6 (was row 14)
$ws.Range("D6").Value = 2.55
$ws.Range("E6").Value = 4.35
$ws.Range("F6").Value = 1.19

# Row 7 (was row 15)
$ws.Range("D7").Value = 0.85
$ws.Range("E7").Value = 0.85
$ws.Range("F7").Value = 0.83

# Row 8 (was row 16)
$ws.Range("D8").Value = 0.67
$ws.Range("E8").Value = 0.67

# Row 9 (was row 17)
$ws.Range("E9").Value = 0.75
$ws.Range("F9").Value = 0.64

# Row 10 (was row 18)
$ws.Range("D10").Value = 0.62
$ws.Range("E10").Value = 0.6
$ws.Range("F10").Value = 0.67

# Print area shrinks to match the new table extent.
$ws.PageSetup.PrintArea = "B2:H10"

# Make this the active sheet/selection (was "Gender & Race" before).
$ws.Activate()
$ws.Range("A5").Select()

$wb.Save()
